# Fruta / hortaliza, semanal
# Inserts two new daily price records (rows 92 and 93) into the
# "Fruta, Agrícola del Norte S.A. de Arica - Plátano" sheet, pushing the
# existing rows 92:171 down to 94:173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 92 (existing row 92 and everything below it
# shifts down by two rows; the new rows inherit formatting from the row
# that follows them, same as Excel's normal "insert row" behaviour).
$ws.Rows.Item(92).Resize(2).Insert()

# --- New row 92 ---------------------------------------------------------
$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(92, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(92, 4).Value = 44554
$ws.Cells.Item(92, 5).Value = 15
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100108
$ws.Cells.Item(92, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(92, 9).Value = 100108006
$ws.Cells.Item(92, 10).Value = "Plátano"
$ws.Cells.Item(92, 11).Value = "Sin especificar"
$ws.Cells.Item(92, 12).Value = "Maduro"
$ws.Cells.Item(92, 13).Value = 60
$ws.Cells.Item(92, 14).Value = 13000
$ws.Cells.Item(92, 15).Value = 13000
$ws.Cells.Item(92, 16).Value = 13000
$ws.Cells.Item(92, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(92, 18).Value = "Bolivia"
$ws.Cells.Item(92, 19).Value = 650
$ws.Cells.Item(92, 20).Value = 20

# --- New row 93 ---------------------------------------------------------
$ws.Cells.Item(93, 1).Value = 1
$ws.Cells.Item(93, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(93, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(93, 4).Value = 44554
$ws.Cells.Item(93, 5).Value = 15
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100108
$ws.Cells.Item(93, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(93, 9).Value = 100108006
$ws.Cells.Item(93, 10).Value = "Plátano"
$ws.Cells.Item(93, 11).Value = "Sin especificar"
$ws.Cells.Item(93, 12).Value = "Pintón"
$ws.Cells.Item(93, 13).Value = 60
$ws.Cells.Item(93, 14).Value = 14000
$ws.Cells.Item(93, 15).Value = 14000
$ws.Cells.Item(93, 16).Value = 14000
$ws.Cells.Item(93, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(93, 18).Value = "Bolivia"
$ws.Cells.Item(93, 19).Value = 700
$ws.Cells.Item(93, 20).Value = 20
